$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plot2")
$ws.Range("H12").Value = 140
